# Weekly update: insert the newest week's price record at row 20 of the
# "Hortaliza, Macroferia Regional de Talca - Arveja Verde" sheet, pushing
# the existing rows 20-52 down to 21-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 20 (shifts rows 20..52 -> 21..53)
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the latest market entry.
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Macroferia Regional de Talca"
$ws.Range("C20").Value = "Maule"
$ws.Range("D20").Value = 44498
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = 100112022
$ws.Range("G20").Value = "Arveja Verde"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 15000
$ws.Range("N20").Value = "$/saco 25 kilos"
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 600
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"
